$wb = $excel.ActiveWorkbook

# Update "展览" sheet (F2: 132 -> 134, F3: 82 -> 83)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 134
$ws1.Range("F3").Value = 83

# Update "全部类型" sheet (F2: 132 -> 134, F3: 82 -> 83)
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value = 134
$ws2.Range("F3").Value = 83
